$d = $word.ActiveDocument

# --- Helper: create a temporary plain-text donor run at the end of the document ---
function New-PlainDonor([string]$text) {
    $lastParaStart = $d.Paragraphs.Item($d.Paragraphs.Count).Range.Start
    $insertionPoint = $d.Range($lastParaStart, $lastParaStart)
    $insertionPoint.InsertBefore($text)
    return $d.Range($lastParaStart, $lastParaStart + $text.Length)
}

# ================= Hunk 1: split "{m" into "{" and "m" =================
$donor1 = New-PlainDonor("{")
$run1 = $d.Range(36, 37)
$run1.Bold = $true
$run1.Bold = $false
$run1.FormattedText = $donor1.FormattedText
$donor1.Delete()

# ================= Hunk 2: split ")}" into ")" and "}" =================
# Use an already-plain run ":" as formatting donor (no text fix needed after,
# since we overwrite its text separately to avoid retrigger of run-merging)
$plainDonorRange = $d.Range(38, 39)   # the ":" run - already plain, no rPr
$run2 = $d.Range(59, 60)              # the "}" character
$run2.FormattedText = $plainDonorRange.FormattedText
$run2.Text = "}"
